# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" for files that are still mid-flight
# (Ready for handoff / Handback transform failed) so the per-locale sheets
# and the Overview roll-up reflect the newest handoff run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# zh-cn: column E is "Latest Handoff Datetime"
$zhcn.Range("E7").Value = "2016-03-19 14:26:31"
$zhcn.Range("E10").Value = "2016-03-19 14:26:31"
$zhcn.Range("E11").Value = "2016-03-19 14:26:31"
$zhcn.Range("E12").Value = "2016-03-19 14:26:31"
$zhcn.Range("E13").Value = "2016-03-19 14:26:31"
$zhcn.Range("E14").Value = "2016-03-19 14:26:31"
$zhcn.Range("E15").Value = "2016-03-19 14:26:31"
$zhcn.Range("E16").Value = "2016-03-19 14:26:31"

# de-de: column E is "Latest Handoff Datetime"
$dede.Range("E7").Value = "2016-03-19 14:26:37"
$dede.Range("E10").Value = "2016-03-19 14:26:37"
$dede.Range("E11").Value = "2016-03-19 14:26:37"
$dede.Range("E12").Value = "2016-03-19 14:26:37"
$dede.Range("E13").Value = "2016-03-19 14:26:37"
$dede.Range("E14").Value = "2016-03-19 14:26:37"
$dede.Range("E15").Value = "2016-03-19 14:26:37"
$dede.Range("E16").Value = "2016-03-19 14:26:37"

# Overview: column D is "Latest Handoff Date", rolled up as the max of the
# per-locale handoff datetimes. Every row that previously shared the old
# "Ready for handoff" timestamp moves forward to the new max.
$overview.Range("D7").Value = "2016-26-19 14:26:37"
$overview.Range("D10").Value = "2016-26-19 14:26:37"
$overview.Range("D11").Value = "2016-26-19 14:26:37"
$overview.Range("D12").Value = "2016-26-19 14:26:37"
$overview.Range("D13").Value = "2016-26-19 14:26:37"
$overview.Range("D14").Value = "2016-26-19 14:26:37"
$overview.Range("D15").Value = "2016-26-19 14:26:37"
$overview.Range("D16").Value = "2016-26-19 14:26:37"
